$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was added to the table. It was inserted right
# after the existing row for "Femacal de La Calera" (row 161), pushing
# every subsequent data row down by one (old row 162 -> new row 163,
# ..., old row 240 -> new row 241).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(162, 1).Value = 3
$ws.Cells.Item(162, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(162, 3).Value = "Coquimbo"
$ws.Cells.Item(162, 4).Value = 44523
$ws.Cells.Item(162, 5).Value = 5
$ws.Cells.Item(162, 6).Value = 100114013
$ws.Cells.Item(162, 7).Value = "Zanahoria"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 448
$ws.Cells.Item(162, 11).Value = 6500
$ws.Cells.Item(162, 12).Value = 7000
$ws.Cells.Item(162, 13).Value = 6701
$ws.Cells.Item(162, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(162, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(162, 16).Value = 335
$ws.Cells.Item(162, 17).Value = 20
$ws.Cells.Item(162, 18).Value = "Hortaliza"
